# Implement "slightly buggy magic projectiles": add a "y" marker to column R
# for most (but not all) rows that already have the pattern filled across
# A:Q, then move the active selection to R20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(2,3,4,5,6,7,9,10,14,15)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 18).Value = "y"
}

$ws.Range("R20").Select()
